# Re-sort the "Tesis publicadas" query table by "Año" (column A) in
# descending order (it was previously sorted ascending), then leave the
# sheet/table as the active selection -- mirroring the manual "Data > Sort
# Z to A" action a user would perform in the Excel UI on that table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tesis publicadas")

# The sheet has a single ListObject (Excel Table) bound to the query table
# "Tesis_publicadas__2" covering A1:F45.
$lo = $ws.ListObjects.Item(1)

# Sort by the first column ("Año") descending, keeping the header row fixed.
$sortColumn = $lo.ListColumns.Item(1).Range
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($sortColumn, 0, 2)
$lo.Sort.Header = 1
$lo.Sort.Apply()

# Make this the active sheet/selection, with the whole table selected --
# matching the state Excel leaves behind right after sorting a table from
# the UI.
$ws.Activate()
$lo.Range.Select()
